$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1000910.2
$ws.Range("I43").Value = 1050
$ws.Range("J43").Value = 1667483.6
$ws.Range("K43").Value = 1050
$ws.Range("L43").Value = 1667483.6
$ws.Range("M43").Value = -981
$ws.Range("N43").Value = -1667621.6
$ws.Range("H129").Value = 1379.9642
$ws.Range("I129").Value = 496.33334
$ws.Range("J129").Value = 1486
$ws.Range("K129").Value = 1489.00002
$ws.Range("L129").Value = 4458
$ws.Range("M129").Value = 3510.99998
$ws.Range("N129").Value = -14458
$ws.Range("H136").Value = 38800
$ws.Range("J136").Value = 38800
$ws.Range("L136").Value = 38800
$ws.Range("N136").Value = -49000
$ws.Range("H137").Value = 1309.1111
$ws.Range("I137").Value = 972
$ws.Range("K137").Value = 2916
$ws.Range("M137").Value = -366
$ws.Range("H139").Value = 39113.332
$ws.Range("J139").Value = 39113.332
$ws.Range("L139").Value = 39113.332
$ws.Range("N139").Value = -49393.332

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5844.085
$ws.Range("I32").Value = 5240.6313
$ws.Range("J32").Value = 8392
$ws.Range("K32").Value = 5240.6313
$ws.Range("L32").Value = 8392
$ws.Range("M32").Value = -4953.6313
$ws.Range("N32").Value = -8966
$ws.Range("H37").Value = 11105.25
$ws.Range("I37").Value = 2739.75
$ws.Range("J37").Value = 15288
$ws.Range("K37").Value = 2739.75
$ws.Range("L37").Value = 15288
$ws.Range("M37").Value = -2466.75
$ws.Range("N37").Value = -15834
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 166676110
$ws.Range("I20").Value = 25403.5
$ws.Range("J20").Value = 250001470
$ws.Range("K20").Value = 25403.5
$ws.Range("L20").Value = 250001470
$ws.Range("M20").Value = -25156.5
$ws.Range("N20").Value = -250001964
$ws.Range("H86").Value = 2294.3
$ws.Range("I86").Value = 2030.375
$ws.Range("J86").Value = 3350
$ws.Range("K86").Value = 2030.375
$ws.Range("L86").Value = 3350
$ws.Range("M86").Value = -907.375
$ws.Range("N86").Value = -5596
$ws.Range("H89").Value = 2294.3
$ws.Range("I89").Value = 2030.375
$ws.Range("J89").Value = 3350
$ws.Range("K89").Value = 10151.875
$ws.Range("L89").Value = 16750
$ws.Range("M89").Value = -4535.875
$ws.Range("N89").Value = -27982
$ws.Range("H134").Value = 1595.6296
$ws.Range("I134").Value = 1372.8182
$ws.Range("K134").Value = 4118.4546
$ws.Range("M134").Value = -1583.4546
$ws.Range("H135").Value = 49639.5
$ws.Range("J135").Value = 49639.5
$ws.Range("L135").Value = 49639.5
$ws.Range("N135").Value = -59779.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3151.45
$ws.Range("I31").Value = 3531.1177
$ws.Range("K31").Value = 3531.1177
$ws.Range("M31").Value = -3236.1177
$ws.Range("H34").Value = 3151.45
$ws.Range("I34").Value = 3531.1177
$ws.Range("K34").Value = 3531.1177
$ws.Range("M34").Value = -3329.1177
$ws.Range("H94").Value = 4312.5454
$ws.Range("I94").Value = 4070.6667
$ws.Range("J94").Value = 4403.25
$ws.Range("K94").Value = 4070.6667
$ws.Range("L94").Value = 4403.25
$ws.Range("M94").Value = -3619.6667
$ws.Range("N94").Value = -5305.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4656.9165
$ws.Range("I3").Value = 2121.4285
$ws.Range("J3").Value = 8206.6
$ws.Range("K3").Value = 6364.2855
$ws.Range("L3").Value = 24619.8
$ws.Range("M3").Value = -6252.2855
$ws.Range("N3").Value = -24843.8
$ws.Range("H11").Value = 80054.39999999999
$ws.Range("I11").Value = 133376.67
$ws.Range("J11").Value = 71
$ws.Range("K11").Value = 400130.01
$ws.Range("L11").Value = 213
$ws.Range("M11").Value = -399990.01
$ws.Range("N11").Value = -493
$ws.Range("H131").Value = 11930395
$ws.Range("I131").Value = 55667068
$ws.Range("J131").Value = 2211.0908
$ws.Range("K131").Value = 167001204
$ws.Range("L131").Value = 6633.2724
$ws.Range("M131").Value = -166996164
$ws.Range("N131").Value = -16713.2724

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 13900
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H70").Value = 4873.3335
$ws.Range("I70").Value = 4100
$ws.Range("K70").Value = 4100
$ws.Range("M70").Value = -3830
$ws.Range("H73").Value = 4873.3335
$ws.Range("I73").Value = 4100
$ws.Range("K73").Value = 4100
$ws.Range("M73").Value = -3164
$ws.Range("H123").Value = 25739.5
$ws.Range("J123").Value = 25739.5
$ws.Range("L123").Value = 25739.5
$ws.Range("N123").Value = -30639.5
$ws.Range("H136").Value = 23000.334
$ws.Range("J136").Value = 23000.334
$ws.Range("L136").Value = 69001.00199999999
$ws.Range("N136").Value = -74101.00199999999
$ws.Range("H138").Value = 35800
$ws.Range("J138").Value = 35800
$ws.Range("L138").Value = 35800
$ws.Range("N138").Value = -46080
$ws.Range("H139").Value = 18633.334
$ws.Range("J139").Value = 18633.334
$ws.Range("L139").Value = 18633.334
$ws.Range("N139").Value = -28913.334

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2344
$ws.Range("H55").Value = 345.3125
$ws.Range("I55").Value = 217
$ws.Range("J55").Value = 422.3
$ws.Range("K55").Value = 217
$ws.Range("L55").Value = 422.3
$ws.Range("M55").Value = -44
$ws.Range("N55").Value = -768.3
$ws.Range("H68").Value = 2268.4211
$ws.Range("I68").Value = 2190
$ws.Range("J68").Value = 2355.5557
$ws.Range("K68").Value = 2190
$ws.Range("L68").Value = 2355.5557
$ws.Range("M68").Value = -1441
$ws.Range("N68").Value = -3853.5557
$ws.Range("H71").Value = 2268.4211
$ws.Range("I71").Value = 2190
$ws.Range("J71").Value = 2355.5557
$ws.Range("K71").Value = 10950
$ws.Range("L71").Value = 11777.7785
$ws.Range("M71").Value = -7206
$ws.Range("N71").Value = -19265.7785
$ws.Range("H80").Value = 22000
$ws.Range("J80").Value = 22000
$ws.Range("L80").Value = 22000
$ws.Range("N80").Value = -24246
$ws.Range("H82").Value = 3244.4443
$ws.Range("I82").Value = 3300
$ws.Range("J82").Value = 3200
$ws.Range("K82").Value = 3300
$ws.Range("L82").Value = 3200
$ws.Range("M82").Value = -2939
$ws.Range("N82").Value = -3922
$ws.Range("H83").Value = 22000
$ws.Range("J83").Value = 22000
$ws.Range("L83").Value = 66000
$ws.Range("N83").Value = -77232
$ws.Range("H85").Value = 3244.4443
$ws.Range("I85").Value = 3300
$ws.Range("J85").Value = 3200
$ws.Range("K85").Value = 3300
$ws.Range("L85").Value = 3200
$ws.Range("M85").Value = -2052
$ws.Range("N85").Value = -5696
$ws.Range("H122").Value = 5910.3335
$ws.Range("I122").Value = 7845.25
$ws.Range("K122").Value = 23535.75
$ws.Range("M122").Value = -21085.75
$ws.Range("H134").Value = 20143
$ws.Range("J134").Value = 20143
$ws.Range("L134").Value = 20143
$ws.Range("N134").Value = -30283
$ws.Range("H138").Value = 34283.25
$ws.Range("J138").Value = 34283.25
$ws.Range("L138").Value = 34283.25
$ws.Range("N138").Value = -44563.25

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 5900
$ws.Range("J39").Value = 5900
$ws.Range("L39").Value = 5900
$ws.Range("N39").Value = -6726
$ws.Range("H69").Value = 24165.857
$ws.Range("J69").Value = 24165.857
$ws.Range("L69").Value = 24165.857
$ws.Range("N69").Value = -25663.857
$ws.Range("H72").Value = 24165.857
$ws.Range("J72").Value = 24165.857
$ws.Range("L72").Value = 72497.571
$ws.Range("N72").Value = -79985.571
